$wb = $excel.ActiveWorkbook

# Update the status text everywhere it appears (shared string is reused
# across the Overview, zh-cn and de-de sheets).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Re-fit the "Status"/locale columns now that the text is shorter.
$wsOverview.Range("E1:F1").EntireColumn.AutoFit() | Out-Null
$wsZhCn.Range("C1").EntireColumn.AutoFit() | Out-Null
$wsDeDe.Range("C1").EntireColumn.AutoFit() | Out-Null
